$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Sample_ID value in B2 from 104 to 105
$ws.Range("B2").Value = 105

# Move the active selection to B2
$ws.Range("B2").Select()
